$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45873, 3852.00481440715, 3768.14814771531, 1944, 4263.05795, 93.1333868045068),
    @(45874, 3852.00481441073, 3719.81309022917, 1944, 4263.05795, 91.1194260757686),
    @(45875, 3852.00481440493, 3689.07692385239, 1944, 4263.05795, 89.8387524769776),
    @(45876, 3852.00481440413, 3680.13923553908, 1944, 4263.05795, 89.4663487972897),
    @(45877, 3852.00481440267, 3635.71238373188, 1944, 4263.05795, 87.6152299720507),
    @(45878, 715.68594436081,  1354.63650390687, 1944, 1483.185599, 7.42233993941916),
    @(45879, 634.559148236726, 1290.27669864056, 1944, 1394.016139, 4.40557039182638),
    @(45880, 3735.57396241807, 3533.22307262223, 1944, 4141.097307, 83.1144340501734),
    @(45881, 3735.57396241807, 3558.21970886817, 1944, 4141.097307, 84.1559605604209),
    @(45882, 3735.57396241807, 3538.61803977859, 1944, 4141.097307, 83.3392243483549),
    @(45883, 3735.57396241807, 3509.22700462974, 1944, 4141.097307, 82.1145978838197),
    @(45884, 602.228776212635, 1490.18154331729, 1944, 1350.706272, 12.2774599626939),
    @(45885, 680.522278765008, 1293.79105591117, 1944, 1436.761544, 4.41793004775684),
    @(45886, 602.228776212635, 1261.94360011396, 1944, 1350.706272, 2.76754566255513)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}
